# Applies the OOXML diff to the document:
#  1. "Just open up a browser..." paragraph gets split into 3 runs with
#     a <w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>
#     pair wrapped around "open up".
#  2. "If you see this, it means you have start the influxdb service / server"
#     paragraph gets split into 3 runs with a
#     <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
#     pair wrapped around "influxdb".
#  3. The "C:/influxdb2> influxd" paragraph gets a
#     <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
#     pair wrapped around the bold "influxd" run.
#  4. The "Once you see this ... Now enter http://localhost:8080 ..."
#     paragraph is reflowed into several runs, gains a
#     <w:lastRenderedPageBreak/>, the word " enter" / trailing space get
#     their own runs, a new bold/size-32 space run is inserted before the
#     hyperlink, the hyperlink's visible text becomes bold/size-32 and is
#     shortened to "http://localhost:808", and a new trailing bold/size-32
#     "6" run (also styled as Hyperlink) is appended after the hyperlink.
#
# NOTE: `Range.InsertXML` on a *Paragraph.Range* that happens to be the very
# last paragraph in the document body spuriously appends an extra empty
# paragraph afterwards. Re-wrapping the same [Start,End) bounds through
# `Document.Range(start, end)` before calling InsertXML avoids that quirk
# (and is harmless for non-last paragraphs too), so every edit below goes
# through a small helper that does exactly that.

$d = $word.ActiveDocument

function Set-ParagraphXml($paraRange, [string]$xml) {
    $clone = $d.Range($paraRange.Start, $paraRange.End)
    $clone.InsertXML($xml)
}

function Find-ParagraphByText([string]$exactText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r") -eq $exactText) {
            return $p.Range
        }
    }
    return $null
}

function Find-ParagraphByPrefix([string]$prefixText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefixText)) {
            return $p.Range
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Change 1: "Just open up a browser and type localhost:8086 ..."
# ---------------------------------------------------------------------
$target = Find-ParagraphByText("Just open up a browser and type localhost:8086 in the address bar.")
if ($target -ne $null) {
    $xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        '<w:p>' +
          '<w:r><w:t xml:space="preserve">Just </w:t></w:r>' +
          '<w:proofErr w:type="gramStart"/>' +
          '<w:r><w:t>open up</w:t></w:r>' +
          '<w:proofErr w:type="gramEnd"/>' +
          '<w:r><w:t xml:space="preserve"> a browser and type localhost:8086 in the address bar.</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    Set-ParagraphXml $target $xml1
}

# ---------------------------------------------------------------------
# Change 2: "If you see this, it means you have start the influxdb ..."
# ---------------------------------------------------------------------
$target = Find-ParagraphByText("If you see this, it means you have start the influxdb service / server")
if ($target -ne $null) {
    $xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        '<w:p>' +
          '<w:r><w:t xml:space="preserve">If you see this, it means you have start the </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>influxdb</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> service / server</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    Set-ParagraphXml $target $xml2
}

# ---------------------------------------------------------------------
# Change 3: "C:/influxdb2> influxd" -> wrap bold "influxd" run with
# spellStart/spellEnd proofErr markers.
# ---------------------------------------------------------------------
$target = Find-ParagraphByText("C:/influxdb2> influxd")
if ($target -ne $null) {
    $xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        '<w:p>' +
          '<w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
          '<w:r><w:t xml:space="preserve">C:/influxdb2&gt; </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>influxd</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
        '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    Set-ParagraphXml $target $xml3
}

# ---------------------------------------------------------------------
# Change 4: "Once you see this, it means the services has started.  Now
# enter http://localhost:8080 to the browser address bar. ..."
# ---------------------------------------------------------------------
$target = Find-ParagraphByPrefix("Once you see this, it means the services has started.")
if ($target -ne $null) {
    $xml4 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>' +
        '<w:p>' +
          '<w:r><w:lastRenderedPageBreak/><w:t>Once you see this, it means the services has started.  Now</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> enter</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
          '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
          '<w:hyperlink r:id="rId6" w:history="1">' +
            '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>http://localhost:808</w:t></w:r>' +
          '</w:hyperlink>' +
          '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>6</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> to the browser address bar.  You will now see InfluxDB2 server admin page on your screen</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    Set-ParagraphXml $target $xml4

    # The InsertXML backend silently drops a bare <w:rStyle> inside <w:rPr>,
    # so the "Hyperlink" character style has to be (re)applied afterwards via
    # the dedicated Range.Style property, which writes <w:rStyle> correctly
    # as long as the range does not also span a paragraph mark.
    $target2 = Find-ParagraphByPrefix("Once you see this, it means the services has started.")

    $searchRange = $d.Range($target2.Start, $target2.End)
    $found = $searchRange.Find.Execute("http://localhost:808", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $searchRange.Style = "Hyperlink"

        # The lone trailing "6" run sits right after the hyperlink run; Word
        # inserts a zero-length boundary position between the two runs, so
        # widen by one extra unit to land on the visible "6" character.
        $sixRange = $d.Range($searchRange.End, $searchRange.End + 2)
        if ($sixRange.Text -eq "6") {
            $sixRange.Style = "Hyperlink"
        } else {
            $sixRange = $d.Range($searchRange.End, $searchRange.End + 1)
            if ($sixRange.Text -eq "6") {
                $sixRange.Style = "Hyperlink"
            }
        }
    }
}
